$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.761.40"
$ws.Range("E2").Value = "  +1.04%  "
$ws.Range("D3").Value = "2.210.83"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'292.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.08%  "
$ws.Range("D6").Value = "'86.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.01%  "
$ws.Range("D7").Value = "'0.515"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.98%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("D10").Value = "'30.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.23%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").Value = "'47.42"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("E14").Value = "  +1.48%  "
$ws.Range("D15").Value = "2.555.40"
$ws.Range("E15").Value = "  +0.98%  "
$ws.Range("D16").Value = "'13.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("D17").Value = "2.210.05"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +2.56%  "
$ws.Range("D19").Value = "39.712.61"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").Value = "'11.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.95%  "
$ws.Range("E22").Value = "  +1.32%  "
$ws.Range("D23").Value = "'65.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("D24").Value = "'235.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").Value = "  +2.68%  "
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("D28").Value = "'22.70"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("E29").Value = "  +0.95%  "
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").Value = "'32.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'151.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +2.76%  "
$ws.Range("D35").Value = "'0.0717"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.63%  "
$ws.Range("E36").Value = "  +1.83%  "
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("D39").Value = "'15.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.91%  "
$ws.Range("E40").Value = "  +3.45%  "
$ws.Range("E41").Value = "  +3.27%  "
$ws.Range("D42").Value = "2.063.23"
$ws.Range("E42").Value = "  +8.53%  "
$ws.Range("E43").Value = "  +4.80%  "
$ws.Range("D45").Value = "'0.0268"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'17.84"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +11.88%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +10.45%  "
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("D49").Value = "2.432.11"
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("D50").Value = "'70.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("D51").Value = "'88.77"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.01%  "
